$d = $word.ActiveDocument

$oldText = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna."
$newText = "Lorem ipsum dolor sit amet, consectetur adipiscing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea."

[void]$d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                               $true, 1, $false, $newText, 2)
